$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item('展览')
$ws1.Range("G2").Value = '已售罄'
$ws1.Range("F3").Value = 7282
$ws1.Range("F4").Value = 3508
$ws1.Range("F6").Value = 3843
$ws1.Range("F8").Value = 78
$ws1.Range("F9").Value = 77
$ws1.Range("F11").Value = 147
$ws1.Range("F12").Value = 506
$ws1.Range("F15").Value = 364
$ws1.Range("F19").Value = 4113
$ws1.Range("F21").Value = 409
$ws1.Range("F23").Value = 535
$ws1.Range("F24").Value = 1654
$ws1.Range("F27").Value = 3027
$ws1.Range("F28").Value = 2217
$ws1.Range("F33").Value = 96
$ws1.Range("F36").Value = 4294
$ws1.Range("F37").Value = 473
$ws1.Range("F41").Value = 798
$ws1.Range("F42").Value = 205
$ws1.Range("F47").Value = 603
$ws1.Range("F48").Value = 715

$ws2 = $wb.Worksheets.Item('演出')
$ws2.Range("F2").Value = 250
$ws2.Range("F15").Value = 27
$ws2.Range("F16").Value = 576

$ws4 = $wb.Worksheets.Item('全部类型')
$ws4.Range("F3").Value = 250
$ws4.Range("C4").Value = '北京·ICOS国际动漫节×CGF中国游戏节02'
$ws4.Range("E4").Value = '2024.06.08 09:00-06.09 17:00'
$ws4.Range("F4").Value = 7282
$ws4.Range("G4").Value = 80
$ws4.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=83161'
$ws4.Range("I4").Value = '//i2.hdslb.com/bfs/openplatform/202405/4uZ0MfIQ1717054288812.jpeg'
$ws4.Range("C5").Value = '北京·thebONE游戏动漫节'
$ws4.Range("D5").Value = '小关路39号 北投购物公园'
$ws4.Range("E5").Value = '2024.06.08 10:00-06.10 17:00'
$ws4.Range("F5").Value = 3508
$ws4.Range("G5").Value = 6.6
$ws4.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=83830'
$ws4.Range("I5").Value = '//i0.hdslb.com/bfs/openplatform/202404/PAQ2DFrV1712046388743.jpeg'
$ws4.Range("F6").Value = 3508
$ws4.Range("C7").Value = '北京·亦创·梦次元动漫游戏展1st'
$ws4.Range("D7").Value = '亦庄荣昌东街6号 北京亦创国际会展中心'
$ws4.Range("E7").Value = '2024.06.08 09:30-06.08 17:00'
$ws4.Range("F7").Value = 3843
$ws4.Range("G7").Value = 80
$ws4.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=84015'
$ws4.Range("I7").Value = '//i1.hdslb.com/bfs/openplatform/202404/UfpmzLsm1712649924888.jpeg'
$ws4.Range("C8").Value = '北京·原神·崩坏·星铁互动展区ONLY'
$ws4.Range("F8").Value = 66
$ws4.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=85926'
$ws4.Range("I8").Value = '//i0.hdslb.com/bfs/openplatform/202405/UPxxwIPm1716180827049.png'
$ws4.Range("F9").Value = 78
$ws4.Range("F10").Value = 77
$ws4.Range("F13").Value = 147
$ws4.Range("F14").Value = 506
$ws4.Range("F17").Value = 364
$ws4.Range("F21").Value = 4113
$ws4.Range("F25").Value = 409
$ws4.Range("F27").Value = 535
$ws4.Range("F28").Value = 1654
$ws4.Range("F31").Value = 3027
$ws4.Range("F32").Value = 2217
$ws4.Range("F39").Value = 4294
$ws4.Range("F41").Value = 473
$ws4.Range("F44").Value = 798
$ws4.Range("F45").Value = 205
$ws4.Range("F49").Value = 603
$ws4.Range("F50").Value = 715
